# Update formatFile and add StaffList.txt
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 3 new rows starting at row 55. This pushes the existing rows
# 55 ("Number", ...), 56, 57 down to 58, 59, 60 respectively, leaving
# row 57 blank (matches the target layout which has a gap at row 57).
$ws.Rows("55:57").Insert()

# New row 55: "Staff List" entry (Name/Format/Example + first struct line)
$ws.Cells.Item(55, 1).Value = "Staff List"
$ws.Cells.Item(55, 2).Value = "StaffList.txt"
$ws.Cells.Item(55, 3).Value = "StaffList.txt"
$ws.Cells.Item(55, 4).Value = "L1"
$ws.Cells.Item(55, 5).Value = """Number of staffs"", staffList.size()"
$ws.Cells.Item(55, 6).Value = "input"

# New row 56: continuation of the Staff struct description
$ws.Cells.Item(56, 4).Value = "RL"
$ws.Cells.Item(56, 5).Value = "ID, password"
$ws.Cells.Item(56, 6).Value = "input"

# "Staff" (the struct/type label) is entered last on both rows so the
# shared-string table gets its new entries in the same order the
# original authoring session produced them in.
$ws.Cells.Item(55, 7).Value = "Staff"
$ws.Cells.Item(56, 7).Value = "Staff"

# Match the author's final selection/scroll position (new rows pushed the
# view down towards the bottom of the sheet).
$excel.ActiveWindow.ScrollRow = 49
$ws.Range("C54").Select()

$wb.Save()
